# Updates the feature-importance table on Sheet1 (A2:D47) with the refreshed
# model results ("finished ls model" run): new gvkey ids (col A), feature
# names (col B), importance scores (col C) and model ids (col D) for all
# 46 ranked rows (rows 2-38 replace the previous 37 rows, rows 39-47 are newly
# added), and extends the sheet's used range/dimension accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object 'object[,]' 46,4
$data[0,0] = 0
$data[0,1] = 'at'
$data[0,2] = 0.1493196360321636
$data[0,3] = 2
$data[1,0] = 12
$data[1,1] = 'dltt_std'
$data[1,2] = 0.05875585310888207
$data[1,3] = 1
$data[2,0] = 18
$data[2,1] = 'icapt_std'
$data[2,2] = 0.05417153320514305
$data[2,3] = 1
$data[3,0] = 36
$data[3,1] = 'sstk_std'
$data[3,2] = 0.04948921474904075
$data[3,3] = 2
$data[4,0] = 39
$data[4,1] = 'teq'
$data[4,2] = 0.04179380529615063
$data[4,3] = 1
$data[5,0] = 23
$data[5,1] = 'pi_std'
$data[5,2] = 0.04138685839021414
$data[5,3] = 2
$data[6,0] = 33
$data[6,1] = 'seq'
$data[6,2] = 0.04119882356890355
$data[6,3] = 1
$data[7,0] = 5
$data[7,1] = 'ceqt'
$data[7,2] = 0.03856959946380981
$data[7,3] = 2
$data[8,0] = 30
$data[8,1] = 'rest_sum_diff'
$data[8,2] = 0.02951416639707502
$data[8,3] = 1
$data[9,0] = 43
$data[9,1] = 'xopr'
$data[9,2] = 0.02894008608697533
$data[9,3] = 1
$data[10,0] = 38
$data[10,1] = 'st_per_growth'
$data[10,2] = 0.02851418192653138
$data[10,3] = 2
$data[11,0] = 16
$data[11,1] = 'gvkey'
$data[11,2] = 0.02726563099250992
$data[11,3] = 1
$data[12,0] = 11
$data[12,1] = 'dltr'
$data[12,2] = 0.02694016422376596
$data[12,3] = 1
$data[13,0] = 4
$data[13,1] = 'ceq'
$data[13,2] = 0.0268362451338258
$data[13,3] = 1
$data[14,0] = 34
$data[14,1] = 'spce'
$data[14,2] = 0.02623585112287411
$data[14,3] = 1
$data[15,0] = 37
$data[15,1] = 'st_per_currentToMax'
$data[15,2] = 0.02418726886050731
$data[15,3] = 1
$data[16,0] = 40
$data[16,1] = 'tstk'
$data[16,2] = 0.02283477198487486
$data[16,3] = 1
$data[17,0] = 14
$data[17,1] = 'fopo_std'
$data[17,2] = 0.02204781222967181
$data[17,3] = 1
$data[18,0] = 2
$data[18,1] = 'caps'
$data[18,2] = 0.02155756183950223
$data[18,3] = 2
$data[19,0] = 32
$data[19,1] = 'sec_trt1m_std'
$data[19,2] = 0.02087360387148816
$data[19,3] = 1
$data[20,0] = 17
$data[20,1] = 'icapt'
$data[20,2] = 0.02024700969472336
$data[20,3] = 1
$data[21,0] = 29
$data[21,1] = 'rest_count_of_diffs'
$data[21,2] = 0.01999105343109417
$data[21,3] = 1
$data[22,0] = 19
$data[22,1] = 'invch'
$data[22,2] = 0.01924200727656406
$data[22,3] = 1
$data[23,0] = 21
$data[23,1] = 'lse'
$data[23,2] = 0.01844090659885377
$data[23,3] = 1
$data[24,0] = 31
$data[24,1] = 'revt'
$data[24,2] = 0.01788921331809505
$data[24,3] = 1
$data[25,0] = 6
$data[25,1] = 'ch'
$data[25,2] = 0.0176477985363721
$data[25,3] = 1
$data[26,0] = 35
$data[26,1] = 'sstk'
$data[26,2] = 0.01751182415085334
$data[26,3] = 2
$data[27,0] = 28
$data[27,1] = 'rest_count'
$data[27,2] = 0.01652372201133656
$data[27,3] = 1
$data[28,0] = 15
$data[28,1] = 'gp'
$data[28,2] = 0.01626000356107158
$data[28,3] = 1
$data[29,0] = 44
$data[29,1] = 'xsga'
$data[29,2] = 0.01491812241746441
$data[29,3] = 1
$data[30,0] = 7
$data[30,1] = 'cogs'
$data[30,2] = 0.01471568354104419
$data[30,3] = 2
$data[31,0] = 42
$data[31,1] = 'wcap'
$data[31,2] = 0.014258735152101
$data[31,3] = 1
$data[32,0] = 27
$data[32,1] = 'rect_std'
$data[32,2] = 0.01421443760922896
$data[32,3] = 1
$data[33,0] = 8
$data[33,1] = 'cogs_std'
$data[33,2] = 0.01397913406897142
$data[33,3] = 1
$data[34,0] = 24
$data[34,1] = 'ppegt'
$data[34,2] = 0.01355190928703611
$data[34,3] = 1
$data[35,0] = 20
$data[35,1] = 'lct'
$data[35,2] = 0.01344980308539786
$data[35,3] = 1
$data[36,0] = 13
$data[36,1] = 'dpact'
$data[36,2] = 0.01339111805122542
$data[36,3] = 1
$data[37,0] = 25
$data[37,1] = 'rat_spcsrc'
$data[37,2] = 0.01316162987401302
$data[37,3] = 1
$data[38,0] = 22
$data[38,1] = 'np_std'
$data[38,2] = 0.01202623454613894
$data[38,3] = 1
$data[39,0] = 26
$data[39,1] = 're'
$data[39,2] = 0.01201839670398554
$data[39,3] = 1
$data[40,0] = 41
$data[40,1] = 'tstk_std'
$data[40,2] = 0.01126560191188868
$data[40,3] = 1
$data[41,0] = 1
$data[41,1] = 'auop'
$data[41,2] = 0.0109743809134068
$data[41,3] = 1
$data[42,0] = 9
$data[42,1] = 'dilavx_std'
$data[42,2] = 0.01094962735025543
$data[42,3] = 1
$data[43,0] = 45
$data[43,1] = 'xsga_std'
$data[43,2] = 0.01061918676823677
$data[43,3] = 1
$data[44,0] = 3
$data[44,1] = 'capx'
$data[44,2] = 0.01036808984658394
$data[44,3] = 1
$data[45,0] = 10
$data[45,1] = 'dlc'
$data[45,2] = 0.0103637576833913
$data[45,3] = 1

$ws.Range("A2:D47").Value2 = $data

# Ensure the newly added rows (39-47) inherit the same formatting as existing column-A cells
# (bold, centered, thin border) by copying the format from an existing styled cell.
$ws.Range("A2").Copy()
$ws.Range("A39:A47").PasteSpecial(-4122)
$excel.CutCopyMode = 0

